$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 138514
$ws.Range("D3").Value = 138515
$ws.Range("D4").Value = 138516
